$wb = $excel.ActiveWorkbook
$wsMain = $wb.Worksheets.Item("main")
$wsModel = $wb.Worksheets.Item("model")

# --- main sheet: Price assumption change (C2 "Price") ---
$wsMain.Range("D2").Value = 413
$wsMain.Range("D2").NumberFormat = "#,##0.00"

# --- model sheet: Maturity assumption change (X85 "Maturity") ---
$wsModel.Range("Y85").Value = 0.04
$wsModel.Range("Y85").NumberFormat = "0.0%"
# Discount (X86 "Discount") keeps its value but gets the same new number format
$wsModel.Range("Y86").NumberFormat = "0.0%"
# Price/NPV ratio cell picks up more decimal precision formatting
$wsModel.Range("Y89").NumberFormat = "#,##0.00"

# --- view-state: selection on model sheet moves, active tab moves to main ---
[void]$wsModel.Range("W86").Select()
[void]$wsMain.Activate()
[void]$wsMain.Range("D3").Select()
